$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46, pushing existing rows 46..102 down to 47..103.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with this week's record (a duplicate of
# the Hass/Primera/Peru record, but dated this week).
$ws.Cells.Item(46, 1).Value = 1
$ws.Cells.Item(46, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value = 44741
$ws.Cells.Item(46, 5).Value = 15
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100106
$ws.Cells.Item(46, 8).Value = "Oleaginosos"
$ws.Cells.Item(46, 9).Value = 100106002
$ws.Cells.Item(46, 10).Value = "Palta"
$ws.Cells.Item(46, 11).Value = "Hass"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 400
$ws.Cells.Item(46, 14).Value = 15000
$ws.Cells.Item(46, 15).Value = 16000
$ws.Cells.Item(46, 16).Value = 15500
$ws.Cells.Item(46, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(46, 18).Value = "Perú"
$ws.Cells.Item(46, 19).Value = 1550
$ws.Cells.Item(46, 20).Value = 10
